# Group all top-level shapes on slide 1 into a single new group shape.
# This mirrors the author selecting every shape on the slide (TextBoxes
# "A"/"B", the two annotated picture groups, the two title textboxes,
# the "C" textbox, and the bottom annotated picture group) and pressing
# Ctrl+G to combine them into one group ("Group 1").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$count = $s.Shapes.Count
$indices = 1..$count

$range = $s.Shapes.Range($indices)
$group = $range.Group()

Write-Host "Grouped" $count "shapes into" $group.Name "(id=" $group.Id ")"
Write-Host "Shapes on slide now:" $s.Shapes.Count
